$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 7, shifting rows 7-49 down to 8-50
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "dct:title"
$ws.Range("B7").Value = "NICEST-2 controlled vocabulary of subjects"
$ws.Range("C7:S7").Value = ""
